$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the computed percent_moisture values (P2:P37) with a placeholder
# constant while more work is done on this column.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 16).Value = 0.00001
}

# Give the percent_moisture header (P1) its own font (Calibri) to set it
# apart while the column is still being finalized.
$ws.Range("P1").Font.Name = "Calibri"

# Move the active selection.
$ws.Range("Q14").Select()
